$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Stash the bold-red-centered "BIC (Combined)" highlight that currently
# ---- lives on C2 ("All parameters different") in a scratch cell, because it
# ---- needs to move onto the new row for that same case (now row 3 / C3).
[void]$ws.Range("C2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# ---- Wipe all existing content+formatting in the used area ----
$ws.Range("A1:H20").Clear()

# ---- Header row ----
$ws.Cells.Item(1,1).Value = "Case #"
$ws.Cells.Item(1,2).Value = "Parameter settings (Linear variation)"
$ws.Cells.Item(1,3).Value = "BIC (Combined)"
$ws.Cells.Item(1,4).Value = "Comments"

# ---- Data rows: Case#, Parameter settings, BIC, Comments ----
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "All parameters different - No storage pool"
$ws.Cells.Item(2,3).Value = 1507
$ws.Cells.Item(2,4).Value = "The model needs the storage pool"

$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "All parameters different - time depandant"
$ws.Cells.Item(3,3).Value = 1373
$ws.Cells.Item(3,4).Value = "Linear variation has the best data fit"

$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "All parameters constant - time depandant"
$ws.Cells.Item(4,3).Value = 3240
$ws.Cells.Item(4,4).Value = "Doesn’t fit the data well enough"

$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "k and Y constant - time depandant"
$ws.Cells.Item(5,3).Value = 1420

$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "Allocations constant - time depandant"
$ws.Cells.Item(6,3).Value = 1392
$ws.Cells.Item(6,4).Value = "Works equally well as case 2, except the roots "

$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "Turnovers constant - time depandant"
$ws.Cells.Item(7,3).Value = 1416

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "All parameters different - temperature depandant"
$ws.Cells.Item(8,3).Value = 1433
$ws.Cells.Item(8,4).Value = "Linear variation has the best data fit"

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "All parameters different - plant size (height) depandant"
$ws.Cells.Item(9,3).Value = 1514
$ws.Cells.Item(9,4).Value = "Linear variation has the best data fit"

# ---- Restore the bold-red-centered format onto the new C3 (BIC for case 2) ----
[void]$ws.Range("Z1").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# ---- Column widths (bestFit); this runtime stores ColumnWidth + 5/6 as the
# ---- raw OOXML <col width>, so compensate to land as close as possible on
# ---- the target stored widths (6.33203125 / 45.83203125 / 13.6640625 / 38.83203125)
$ws.Columns.Item(1).ColumnWidth = 5.498697916666667
$ws.Columns.Item(2).ColumnWidth = 44.998697916666664
$ws.Columns.Item(3).ColumnWidth = 12.830729166666666
$ws.Columns.Item(4).ColumnWidth = 37.998697916666664

# ---- Selection ----
[void]$ws.Range("E3").Select()

Write-Output "done"
